$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-60 (values shifted / changed per new weekly data) ---
$ws.Range("D8").Value = 44329
$ws.Range("J8").Value = 900
$ws.Range("K8").Value = 350
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 1633
$ws.Range("P8").Value = 33
$ws.Range("D9").Value = 44364
$ws.Range("K9").Value = 8500
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 8750
$ws.Range("N9").Value = "`$/caja 50 unidades"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("Q9").Value = 50
$ws.Range("D10").Value = 44258
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 10500
$ws.Range("N10").Value = "`$/caja 60 unidades"
$ws.Range("P10").Value = 175
$ws.Range("Q10").Value = 60
$ws.Range("D11").Value = 44159
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 6500
$ws.Range("O11").Value = "Región de O'Higgins"
$ws.Range("P11").Value = 130
$ws.Range("D12").Value = 44435
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 12500
$ws.Range("N12").Value = "`$/caja 50 unidades"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 250
$ws.Range("Q12").Value = 50
$ws.Range("D13").Value = 44195
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 10500
$ws.Range("P13").Value = 175
$ws.Range("D14").Value = 44230
$ws.Range("J14").Value = 150
$ws.Range("M14").Value = 9333
$ws.Range("P14").Value = 156
$ws.Range("D15").Value = 44272
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 9500
$ws.Range("O15").Value = "Región de O'Higgins"
$ws.Range("P15").Value = 158
$ws.Range("D16").Value = 44293
$ws.Range("K16").Value = 8000
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 8500
$ws.Range("N16").Value = "`$/caja 60 unidades"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 142
$ws.Range("Q16").Value = 60
$ws.Range("D17").Value = 44355
$ws.Range("K17").Value = 9000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 9500
$ws.Range("N17").Value = "`$/caja 50 unidades"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 190
$ws.Range("Q17").Value = 50
$ws.Range("D18").Value = 44238
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 11000
$ws.Range("M18").Value = 10500
$ws.Range("N18").Value = "`$/caja 60 unidades"
$ws.Range("O18").Value = "Región de O'Higgins"
$ws.Range("P18").Value = 175
$ws.Range("Q18").Value = 60
$ws.Range("D19").Value = 44299
$ws.Range("K19").Value = 7000
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = 7500
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 150
$ws.Range("D20").Value = 44320
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 8500
$ws.Range("P20").Value = 170
$ws.Range("D21").Value = 44385
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 9500
$ws.Range("N21").Value = "`$/caja 50 unidades"
$ws.Range("O21").Value = "Región de Arica y Parinacota"
$ws.Range("P21").Value = 190
$ws.Range("Q21").Value = 50
$ws.Range("D22").Value = 44253
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 8500
$ws.Range("N22").Value = "`$/caja 60 unidades"
$ws.Range("O22").Value = "Región de O'Higgins"
$ws.Range("P22").Value = 142
$ws.Range("Q22").Value = 60
$ws.Range("D23").Value = 44334
$ws.Range("K23").Value = 11000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 11500
$ws.Range("N23").Value = "`$/caja 50 unidades"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 230
$ws.Range("Q23").Value = 50
$ws.Range("D24").Value = 44281
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 9500
$ws.Range("N24").Value = "`$/caja 60 unidades"
$ws.Range("O24").Value = "Región de O'Higgins"
$ws.Range("P24").Value = 158
$ws.Range("Q24").Value = 60
$ws.Range("D25").Value = 44434
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 12500
$ws.Range("P25").Value = 250
$ws.Range("D26").Value = 44327
$ws.Range("D27").Value = 44383
$ws.Range("N27").Value = "`$/caja 50 unidades"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 190
$ws.Range("Q27").Value = 50
$ws.Range("D28").Value = 44223
$ws.Range("K28").Value = 9000
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 9500
$ws.Range("N28").Value = "`$/caja 60 unidades"
$ws.Range("O28").Value = "Región de O'Higgins"
$ws.Range("P28").Value = 158
$ws.Range("Q28").Value = 60
$ws.Range("D29").Value = 44441
$ws.Range("K29").Value = 14000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 14500
$ws.Range("P29").Value = 290
$ws.Range("D30").Value = 44336
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = 11000
$ws.Range("P30").Value = 220
$ws.Range("D31").Value = 44341
$ws.Range("K31").Value = 9000
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = 9500
$ws.Range("P31").Value = 190
$ws.Range("D32").Value = 44453
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = 15500
$ws.Range("N32").Value = "`$/caja 50 unidades"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 310
$ws.Range("Q32").Value = 50
$ws.Range("D33").Value = 44208
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11500
$ws.Range("N33").Value = "`$/caja 60 unidades"
$ws.Range("O33").Value = "Región de O'Higgins"
$ws.Range("P33").Value = 192
$ws.Range("Q33").Value = 60
$ws.Range("D34").Value = 44420
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 9500
$ws.Range("P34").Value = 190
$ws.Range("D35").Value = 44370
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 11000
$ws.Range("M35").Value = 10500
$ws.Range("N35").Value = "`$/caja 50 unidades"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 210
$ws.Range("Q35").Value = 50
$ws.Range("D36").Value = 44237
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 8000
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = 8500
$ws.Range("P36").Value = 142
$ws.Range("D37").Value = 44285
$ws.Range("J37").Value = 100
$ws.Range("O37").Value = "Región de O'Higgins"
$ws.Range("D38").Value = 44217
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 9000
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = 9500
$ws.Range("N38").Value = "`$/caja 60 unidades"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 158
$ws.Range("Q38").Value = 60
$ws.Range("D39").Value = 44455
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 16000
$ws.Range("L39").Value = 17000
$ws.Range("M39").Value = 16500
$ws.Range("N39").Value = "`$/caja 50 unidades"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 330
$ws.Range("Q39").Value = 50
$ws.Range("D40").Value = 44427
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 13000
$ws.Range("M40").Value = 12500
$ws.Range("P40").Value = 250
$ws.Range("D41").Value = 44265
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 7000
$ws.Range("L41").Value = 8000
$ws.Range("M41").Value = 7500
$ws.Range("N41").Value = "`$/caja 60 unidades"
$ws.Range("O41").Value = "Región de O'Higgins"
$ws.Range("P41").Value = 125
$ws.Range("Q41").Value = 60
$ws.Range("D42").Value = 44343
$ws.Range("N42").Value = "`$/caja 50 unidades"
$ws.Range("O42").Value = "Región de Arica y Parinacota"
$ws.Range("P42").Value = 190
$ws.Range("Q42").Value = 50
$ws.Range("D43").Value = 44447
$ws.Range("J43").Value = 100
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 16000
$ws.Range("M43").Value = 15500
$ws.Range("N43").Value = "`$/caja 50 unidades"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value = 310
$ws.Range("Q43").Value = 50
$ws.Range("D44").Value = 44260
$ws.Range("N44").Value = "`$/caja 60 unidades"
$ws.Range("O44").Value = "Región de O'Higgins"
$ws.Range("P44").Value = 158
$ws.Range("Q44").Value = 60
$ws.Range("D45").Value = 44187
$ws.Range("K45").Value = 8000
$ws.Range("L45").Value = 9000
$ws.Range("M45").Value = 8500
$ws.Range("P45").Value = 142
$ws.Range("D46").Value = 44390
$ws.Range("K46").Value = 9000
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = 9500
$ws.Range("N46").Value = "`$/caja 50 unidades"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 190
$ws.Range("Q46").Value = 50
$ws.Range("D47").Value = 44251
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 10000
$ws.Range("M47").Value = 9500
$ws.Range("N47").Value = "`$/caja 60 unidades"
$ws.Range("P47").Value = 158
$ws.Range("Q47").Value = 60
$ws.Range("D48").Value = 44243
$ws.Range("J48").Value = 100
$ws.Range("K48").Value = 10000
$ws.Range("L48").Value = 11000
$ws.Range("M48").Value = 10500
$ws.Range("N48").Value = "`$/caja 60 unidades"
$ws.Range("P48").Value = 175
$ws.Range("Q48").Value = 60
$ws.Range("D49").Value = 44166
$ws.Range("J49").Value = 200
$ws.Range("K49").Value = 6000
$ws.Range("L49").Value = 7000
$ws.Range("M49").Value = 6500
$ws.Range("O49").Value = "Región de O'Higgins"
$ws.Range("P49").Value = 130
$ws.Range("D50").Value = 44168
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 6000
$ws.Range("L50").Value = 7000
$ws.Range("M50").Value = 6500
$ws.Range("O50").Value = "Región de O'Higgins"
$ws.Range("P50").Value = 130
$ws.Range("D51").Value = 44397
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = 8500
$ws.Range("P51").Value = 170
$ws.Range("D52").Value = 44363
$ws.Range("K52").Value = 9000
$ws.Range("L52").Value = 10000
$ws.Range("M52").Value = 9500
$ws.Range("N52").Value = "`$/caja 50 unidades"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value = 190
$ws.Range("Q52").Value = 50
$ws.Range("D53").Value = 44372
$ws.Range("N53").Value = "`$/caja 50 unidades"
$ws.Range("O53").Value = "Región de Arica y Parinacota"
$ws.Range("P53").Value = 190
$ws.Range("Q53").Value = 50
$ws.Range("D54").Value = 44306
$ws.Range("N54").Value = "`$/caja 60 unidades"
$ws.Range("O54").Value = "Región de O'Higgins"
$ws.Range("P54").Value = 142
$ws.Range("Q54").Value = 60
$ws.Range("D55").Value = 44215
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = 9500
$ws.Range("P55").Value = 158
$ws.Range("D56").Value = 44357
$ws.Range("J56").Value = 100
$ws.Range("K56").Value = 8000
$ws.Range("L56").Value = 9000
$ws.Range("M56").Value = 8500
$ws.Range("N56").Value = "`$/caja 50 unidades"
$ws.Range("O56").Value = "Región de Arica y Parinacota"
$ws.Range("P56").Value = 170
$ws.Range("Q56").Value = 50
$ws.Range("D57").Value = 44203
$ws.Range("J57").Value = 200
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 11000
$ws.Range("M57").Value = 10500
$ws.Range("N57").Value = "`$/caja 60 unidades"
$ws.Range("O57").Value = "Región de O'Higgins"
$ws.Range("P57").Value = 175
$ws.Range("Q57").Value = 60
$ws.Range("D58").Value = 44162
$ws.Range("K58").Value = 6000
$ws.Range("L58").Value = 6500
$ws.Range("M58").Value = 6250
$ws.Range("P58").Value = 104
$ws.Range("D59").Value = 44358
$ws.Range("K59").Value = 9000
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = 9500
$ws.Range("N59").Value = "`$/caja 50 unidades"
$ws.Range("P59").Value = 190
$ws.Range("Q59").Value = 50
$ws.Range("D60").Value = 44211
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 11000
$ws.Range("L60").Value = 12000
$ws.Range("M60").Value = 11500
$ws.Range("N60").Value = "`$/caja 60 unidades"
$ws.Range("O60").Value = "Región de O'Higgins"
$ws.Range("P60").Value = 192
$ws.Range("Q60").Value = 60

# --- Add new rows 61 and 62 (new weekly observations) ---
# Row 61
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44425
$ws.Range("D61").NumberFormat = $ws.Range("D2").NumberFormat()
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112032
$ws.Range("G61").Value = "Zapallo italiano"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 11000
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = 11500
$ws.Range("N61").Value = "`$/caja 60 unidades"
$ws.Range("O61").Value = "Región de Arica y Parinacota"
$ws.Range("P61").Value = 192
$ws.Range("Q61").Value = 60
$ws.Range("R61").Value = "Hortaliza"

# Row 62
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 44323
$ws.Range("D62").NumberFormat = $ws.Range("D2").NumberFormat()
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 100112032
$ws.Range("G62").Value = "Zapallo italiano"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 100
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = 9500
$ws.Range("N62").Value = "`$/caja 50 unidades"
$ws.Range("O62").Value = "Región de Arica y Parinacota"
$ws.Range("P62").Value = 190
$ws.Range("Q62").Value = 50
$ws.Range("R62").Value = "Hortaliza"

